# Regenerate save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
# This updates the "K" column (column G) values on Sheet1 to the newly
# recalculated strikeout counts.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 1
    3  = 0
    4  = 1
    5  = 0
    7  = 0
    8  = 1
    9  = 0
    10 = 1
    11 = 1
    12 = 0
    13 = 0
    14 = 1
    15 = 1
    16 = 0
    17 = 2
    18 = 1
    19 = 1
    20 = 1
    21 = 0
    23 = 2
    24 = 0
    25 = 1
    26 = 2
    27 = 2
    28 = 3
    29 = 1
    30 = 2
    31 = 1
    32 = 3
    33 = 1
    34 = 3
    35 = 0
    36 = 0
    37 = 0
    38 = 0
    39 = 0
    40 = 1
    41 = 2
    42 = 3
    43 = 1
    44 = 2
    45 = 2
    46 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
